# Auto update stock data
# Updates the Date_1 (col A) and EBITDA (col B) columns for each ticker's
# most recent row, rolling the data date from 2025/11/03 to 2025/11/04 and
# refreshing the EBITDA figure. Values are written as text (matching the
# existing inline-string cell type) via the leading-apostrophe text-entry
# trick, then ClearFormats() strips the resulting quote-prefix style so the
# cell format stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# row -> new EBITDA value (row 38 and 56 only change the date, EBITDA stays)
$updates = @{
    2  = "4.53"
    8  = "7.57"
    14 = "2.69"
    20 = "12.45"
    26 = "9.81"
    32 = "25.03"
    38 = $null
    44 = "12.04"
    50 = "11.30"
    56 = $null
    62 = "11.94"
    68 = "13.28"
    74 = "15.53"
}

foreach ($row in $updates.Keys | Sort-Object) {
    Set-TextValue $ws.Range("A$row") "2025/11/04"

    $newB = $updates[$row]
    if ($null -ne $newB) {
        Set-TextValue $ws.Range("B$row") $newB
    }
}
